# Daily attendance processing - 2025-12-18 14:59:29
#
# Normalises the "Recorded By" column (G): wherever the literal
# "System" entry is present but not already first in the comma-
# separated list of recorders, move it to the front while leaving the
# relative order of the remaining entries untouched.
#
# Note: comparisons use the .Equals() string method (exact, ordinal,
# case-sensitive) rather than the -eq/-ceq operators, because a
# lower-case "system" value must NOT be treated as a match for the
# capitalised "System" token.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = @($text -split ", ")

    # Find the first part that is exactly "System" (case-sensitive).
    $systemIndex = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i].Equals("System")) {
            $systemIndex = $i
            break
        }
    }

    # Only rewrite when "System" exists and is not already first.
    if ($systemIndex -gt 0) {
        $newParts = @("System")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $systemIndex) {
                $newParts += $parts[$i]
            }
        }
        $cell.Value = $newParts -join ", "
        $changed++
    }
}

Write-Output ("Recorded By column normalised; rows changed: " + $changed)
